$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 entirely (data for these sending/target cluster
# combinations is no longer present after the TPM update).
# Delete bottom-up so row indices of remaining rows don't shift unexpectedly.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Update row 2 with the refreshed TPM-derived values
# (columns K-P are unchanged by this edit and are intentionally left alone)
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Tg"
$ws.Range("C2").Value = "Asgr1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6951913333333333
$ws.Range("H2").Value = 2.085574
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("Q2").Value = 0.7219962890335555
$ws.Range("R2").Value = 6.497966601302
$ws.Range("S2").Value = 0.961760453297246
$ws.Range("T2").Value = 0.961760453297246

# Update row 3 with the refreshed TPM-derived values
# (columns K-P are unchanged by this edit and are intentionally left alone)
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Tg"
$ws.Range("C3").Value = "Asgr1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6951913333333333
$ws.Range("H3").Value = 2.085574
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("Q3").Value = 0.02870653572733333
$ws.Range("R3").Value = 0.258358821546
$ws.Range("S3").Value = 0.03823954670275396
$ws.Range("T3").Value = 0.03823954670275396
